# Atualizacao de bases das ligas, do dia: 14-05-2024 as 20:19
# Append 5 new finished-match rows (166-170) to the bottom of the
# "Azerbaijan Premier League" data table, matching the formatting of the
# existing rows (column A: bold/bordered/centered "id" style; column D:
# "YYYY-MM-DD HH:MM:SS" date style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ row=166; A=164; B=7133550; C="Azerbaijan Premier League"; D=45422.5; E="Neftchi Baku"; F="FK Kapaz"; G=5; H=1; I="H"; J=1.615; K=3.4; L=5; M=1.222; N=5; O=9.5; P=-1.75; Q=1.75; R=1.95; S=3.25; T=1.9; U=1.9; V=0.222; W=-1; X=-1; Y=0.75; Z=-1; AA=0.8999999999999999; AB=-1 },
    @{ row=167; A=165; B=7123423; C="Azerbaijan Premier League"; D=45423.41666666666; E="Araz FK"; F="PFK Turan Tovuz"; G=0; H=1; I="A"; J=2.625; K=3; L=2.5; M=2.7; N=3; O=2.45; P=0; Q=2; R=1.8; S=2.25; T=1.9; U=1.9; V=-1; W=-1; X=1.45; Y=-1; Z=0.8; AA=-1; AB=0.8999999999999999 },
    @{ row=168; A=166; B=7128941; C="Azerbaijan Premier League"; D=45423.52083333334; E="Zira IK"; F="FK Qarabag"; G=0; H=1; I="A"; J=3.2; K=3.6; L=1.909; M=4.2; N=3.8; O=1.615; P=0.75; Q=2; R=1.8; S=2.5; T=1.875; U=1.925; V=-1; W=-1; X=0.615; Y=-0.5; Z=0.4; AA=-1; AB=0.925 },
    @{ row=169; A=167; B=7128942; C="Azerbaijan Premier League"; D=45424.41666666666; E="FK Gabala"; F="Sabah"; G=2; H=0; I="H"; J=7; K=6; L=1.25; M=3.4; N=4.5; O=1.615; P=0.75; Q=1.975; R=1.825; S=2.75; T=1.85; U=1.95; V=2.4; W=-1; X=-1; Y=0.9750000000000001; Z=-1; AA=-1; AB=0.95 },
    @{ row=170; A=168; B=7123424; C="Azerbaijan Premier League"; D=45424.52083333334; E="FK Sumqayit"; F="Sabail FC"; G=2; H=1; I="H"; J=1.8; K=3.8; L=3.4; M=1.85; N=3.75; O=3.3; P=-0.5; Q=1.875; R=1.925; S=2.5; T=1.95; U=1.85; V=0.8500000000000001; W=-1; X=-1; Y=0.875; Z=-1; AA=0.95; AB=-1 }
)

foreach ($r in $newRows) {
    $row = $r.row

    # Column A ("id"): copy the style of the last existing row, then set the
    # value, so the new cell reuses the existing bold/bordered/centered format.
    $ws.Range("A164").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $r.A

    # Column D ("Date"): copy the style of the last existing row, then set the
    # value, so the new cell reuses the existing date number format.
    $ws.Range("D164").Copy()
    $ws.Range("D$row").PasteSpecial(-4122)
    $ws.Range("D$row").Value = $r.D

    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
    $ws.Range("W$row").Value = $r.W
    $ws.Range("X$row").Value = $r.X
    $ws.Range("Y$row").Value = $r.Y
    $ws.Range("Z$row").Value = $r.Z
    $ws.Range("AA$row").Value = $r.AA
    $ws.Range("AB$row").Value = $r.AB
}
